$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple numeric updates (no row reordering) ---

# Row 32: Pakistan
$ws.Range("B32").Value = 12227
$ws.Range("C32").Value = 287
$ws.Range("E32").Value = 9216

# Row 35: Rumania
$ws.Range("E35").Value = 7166
$ws.Range("G35").Value = 12
$ws.Range("H35").Value = 579

# Row 68: Uzbekistan
$ws.Range("D68").Value = 707
$ws.Range("E68").Value = 1121

# --- Senegal overtakes Honduras & Uruguay in the ranking ---
# Old order (rows 102-104): Honduras, Uruguay, Senegal
# New order (rows 102-104): Senegal, Honduras, Uruguay
# Honduras & Uruguay keep their previous stats (shifted down a row);
# Senegal gets freshly updated figures.

$ws.Range("A102").Value = "Senegal"
$ws.Range("B102").Value = 614
$ws.Range("C102").Value = 69
$ws.Range("D102").Value = 276
$ws.Range("E102").Value = 331
$ws.Range("F102").Value = 1
$ws.Range("G102").Value = 0
$ws.Range("H102").Value = 7

$ws.Range("A103").Value = "Honduras"
$ws.Range("B103").Value = 591
$ws.Range("C103").Value = 29
$ws.Range("D103").Value = 58
$ws.Range("E103").Value = 478
$ws.Range("F103").Value = 10
$ws.Range("G103").Value = 8
$ws.Range("H103").Value = 55

$ws.Range("A104").Value = "Uruguay"
$ws.Range("B104").Value = 563
$ws.Range("C104").Value = 0
$ws.Range("D104").Value = 369
$ws.Range("E104").Value = 182
$ws.Range("F104").Value = 9
$ws.Range("G104").Value = 0
$ws.Range("H104").Value = 12

# --- Madagascar overtakes Etiopia in the ranking ---
# Old order (rows 140-141): Etiopia, Madagascar
# New order (rows 140-141): Madagascar, Etiopia
# Etiopia keeps its previous stats (shifted down a row);
# Madagascar gets freshly updated figures.

$ws.Range("A140").Value = "Madagascar"
$ws.Range("B140").Value = 123
$ws.Range("C140").Value = 1
$ws.Range("D140").Value = 62
$ws.Range("E140").Value = 61
$ws.Range("F140").Value = 1
$ws.Range("G140").Value = 0
$ws.Range("H140").Value = 0

$ws.Range("A141").Value = "Etiopia"
$ws.Range("B141").Value = 122
$ws.Range("C141").Value = 5
$ws.Range("D141").Value = 29
$ws.Range("E141").Value = 90
$ws.Range("F141").Value = 0
$ws.Range("G141").Value = 0
$ws.Range("H141").Value = 3

# --- Update the "last refreshed" timestamp string ---
$ws.Range("A1").Value = "Datos actualizados a 25 de Abril de 2020 a las 13:52"
